# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    current "2022-Q3" sheet), seeded from a copy of "2022-Q3" so it keeps
#    the same header / row styling, then overwrite its cell contents with
#    the new quarter's fund-holding data and trim the extra rows.
# 2. Prepend a "2022-Q4" row to the "总计" summary sheet (row 2), pushing
#    the existing quarter rows down by one and adding the trailing
#    "2021-Q3" row that falls out the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q4" sheet
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Copy($templateSheet, $null)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# the template ("2022-Q3") has 13 data rows (rows 2-13); the new
# "2022-Q4" sheet only needs 8 (rows 2-9) -- drop the extra 4 rows.
$newSheet.Rows.Item(10).Resize(4).Delete()

$q4Data = @(
    @("200015", "长城优化升级混合A",     "18.10", "81.69", "3.08", "0.5575", 10),
    @("013274", "长城优化升级混合C",     "7.83",  "81.69", "3.08", "0.2412", 10),
    @("506008", "长城科创两年定开混合A", "3.19",  "77.75", "3.56", "0.1136", 5),
    @("007133", "嘉实长青竞争优势股票A", "0.46",  "91.19", "5.47", "0.0252", 7),
    @("003670", "中融物联网主题灵活配置混合", "0.13", "92.35", "4.41", "0.0057", 1),
    @("501002", "长信价值优选混合",     "0.40",  "81.49", "1.23", "0.0049", 7),
    @("012793", "长城科创两年定开混合C", "0.11",  "77.75", "3.56", "0.0039", 5),
    @("007134", "嘉实长青竞争优势股票C", "0.04",  "91.19", "5.47", "0.0022", 7)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $i + 2
    $row = $q4Data[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: "总计" summary sheet -- insert the 2022-Q4 row at the top of
# the data (row 2), shifting everything else down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 8
$summary.Cells.Item(2, 4).Value = 0.95

# renumber the "序号"-style first column (A) for the rows that shifted
# down, and append the trailing 2021-Q3 row that fell off the bottom.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4

$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(7, 2).Value = "2021-Q3"
$summary.Cells.Item(7, 3).Value = 3
$summary.Cells.Item(7, 4).Value = 0
